$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the "Arquivos reservados" translator column (D) for rows that
#    didn't have it yet, marking them with the reserved-file tag "cfarl".
# ---------------------------------------------------------------------------
for ($r = 648; $r -le 663; $r++) {
    $ws.Cells.Item($r, 4).Value = "cfarl"
}

# ---------------------------------------------------------------------------
# 2) Flip the "traduzido?" column (C) from "Nao" to "SIM" for the batch of
#    rows that got finished.
# ---------------------------------------------------------------------------
$cRows = @(665,666,667,668,669,670,671,672,673,674,675,676,677,678,679,680,681,
           682,683,684,685,688,689,690,691,692,719,720)
foreach ($r in $cRows) {
    $ws.Cells.Item($r, 3).Value = "SIM"
}

# ---------------------------------------------------------------------------
# 3) Conditional-formatting experimentation ("temp cfarl" per the commit):
#    try out a handful of extra "Equal To" highlight rules on column C, then
#    remove them again, leaving only the two real rules this sheet always
#    had (colour-scale + "SIM" highlight on C, not-blank highlight on D).
# ---------------------------------------------------------------------------
$rangeC = $ws.Range("C1:C1048576")

$tmp1 = $rangeC.FormatConditions.Add(1, 3, '="SIM"')
$tmp1.Interior.Color = 13561798

$tmp2 = $rangeC.FormatConditions.Add(1, 3, '="SIM"')
$tmp2.Interior.Color = 13561798

$tmp3 = $rangeC.FormatConditions.Add(1, 3, '="SIM"')
$tmp3.Interior.Color = 13561798
$tmp3.Font.Color = 24832

$tmp4 = $rangeC.FormatConditions.Add(1, 3, '="SIM"')
$tmp4.Interior.Color = 13561798
$tmp4.Font.Color = 24832

$tmp5 = $rangeC.FormatConditions.Add(1, 3, '="SIM"')
$tmp5.Interior.Color = 13561798
$tmp5.Font.Color = 24832

# discard the scratch rules again, from the end so indices stay valid
for ($i = $rangeC.FormatConditions.Count; $i -ge 3; $i--) {
    $rangeC.FormatConditions.Item($i).Delete()
}

# ---------------------------------------------------------------------------
# 4) Leave the view parked on the cell the author was last working on, with
#    the grid scrolled so row 679 is at the top (best-effort: the host may
#    not expose/persist window-scroll state, so don't let that fail the
#    script).
# ---------------------------------------------------------------------------
try {
    $aw = $excel.ActiveWindow
    $aw.ScrollRow = 679
    $aw.ScrollColumn = 1
} catch {
}

$ws.Range("C665").Select() | Out-Null
